$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) to the s_vals sheet, matching the header
# formatting used by the other header cells (e.g. column G "sum").
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Corresponding data value for the new column in row 2.
$ws.Range("H2").Value = 0
